$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B (SR-ID) and C (SIQ-ID) so the order becomes CR-ID, SIQ-ID, SR-ID
# Step 1: stage old column C into helper column Z
$ws.Range("C1:C62").Copy()
$ws.Range("Z1:Z62").Select()
$ws.Paste()

# Step 2: copy old column B into column C
$ws.Range("B1:B62").Copy()
$ws.Range("C1:C62").Select()
$ws.Paste()

# Step 3: copy staged old column C (now in Z) into column B
$ws.Range("Z1:Z62").Copy()
$ws.Range("B1:B62").Select()
$ws.Paste()

# Step 4: cells whose source column was originally blank must be blanked out
# (paste of a blank cell keeps old content, so clear residual values explicitly)
# Rows where old column B was blank -> new column C must be blank
$ws.Range("C8").Value2 = ""
$ws.Range("C13").Value2 = ""
$ws.Range("C19").Value2 = ""
$ws.Range("C22").Value2 = ""
$ws.Range("C28").Value2 = ""
$ws.Range("C36").Value2 = ""
$ws.Range("C44").Value2 = ""
$ws.Range("C49").Value2 = ""
$ws.Range("C56").Value2 = ""

# Rows where old column C was blank -> new column B must be blank
$ws.Range("B6").Value2 = ""
$ws.Range("B8").Value2 = ""
$ws.Range("B11").Value2 = ""
$ws.Range("B13").Value2 = ""
$ws.Range("B17").Value2 = ""
$ws.Range("B19").Value2 = ""
$ws.Range("B20").Value2 = ""
$ws.Range("B21").Value2 = ""
$ws.Range("B22").Value2 = ""
$ws.Range("B28").Value2 = ""
$ws.Range("B35").Value2 = ""
$ws.Range("B36").Value2 = ""
$ws.Range("B43").Value2 = ""
$ws.Range("B44").Value2 = ""
$ws.Range("B47").Value2 = ""
$ws.Range("B48").Value2 = ""
$ws.Range("B49").Value2 = ""
$ws.Range("B53").Value2 = ""
$ws.Range("B55").Value2 = ""
$ws.Range("B56").Value2 = ""
$ws.Range("B60").Value2 = ""
$ws.Range("B62").Value2 = ""

# Step 5: clean up helper column
$ws.Range("Z1:Z62").Clear()

